$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing validation values (rows 2-10)
$ws.Range("B2").Value = 0.5225164411845634
$ws.Range("B3").Value = 0.6164424751215134
$ws.Range("B4").Value = 0.2634196246095795
$ws.Range("B5").Value = 0.05601544662680034
$ws.Range("B6").Value = 0.380197141932949
$ws.Range("B7").Value = 0.8019352917824528
$ws.Range("B8").Value = 0.1637892253675666
$ws.Range("B9").Value = 0.9122409358008438
$ws.Range("B10").Value = 0.4200713990676817

# Add new rows for 6_1 and 6_2
$ws.Range("A11").Value = "6_1"
$ws.Range("B11").Value = 0.9030506874373044
$ws.Range("A12").Value = "6_2"
$ws.Range("B12").Value = 0.3595964676950594
